$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 611, pushing the
# existing rows 611-632 down to 612-633 (all their data stays identical,
# only their row numbers shift).
$ws.Rows.Item(611).Insert()

# Capture the date number format used by the other rows in column D so the
# newly inserted date cell matches the existing formatting/style.
$dateFormat = $ws.Cells.Item(612, 4).NumberFormat

# Populate the new row 611 with the new record's data.
$ws.Cells.Item(611, 1).Value = 10
$ws.Cells.Item(611, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(611, 3).Value = "La Araucanía"
$ws.Cells.Item(611, 4).Value = 45075
$ws.Cells.Item(611, 4).NumberFormat = $dateFormat
$ws.Cells.Item(611, 5).Value = 9
$ws.Cells.Item(611, 6).Value = 100112024
$ws.Cells.Item(611, 7).Value = "Choclo"
$ws.Cells.Item(611, 8).Value = "Dulce o Americano"
$ws.Cells.Item(611, 9).Value = "Primera"
$ws.Cells.Item(611, 10).Value = 185
$ws.Cells.Item(611, 11).Value = 14000
$ws.Cells.Item(611, 12).Value = 15000
$ws.Cells.Item(611, 13).Value = 14324
$ws.Cells.Item(611, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(611, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(611, 16).Value = 205
$ws.Cells.Item(611, 17).Value = 70
$ws.Cells.Item(611, 18).Value = "Hortaliza"
